$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.245.93'
$ws.Range('E2').Value = '  +1.25%  '
$ws.Range('D3').Value = '2.903.60'
$ws.Range('E3').Value = '  +3.93%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'353.76"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.38%  '
$ws.Range('D6').Value = "'113.63"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.19%  '
$ws.Range('D7').Value = "'0.557"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.59%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = "'0.626"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.19%  '
$ws.Range('D10').Value = "'40.05"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.16%  '
$ws.Range('D11').Value = "'0.0865"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.79%  '
$ws.Range('E12').Value = '  +0.62%  '
$ws.Range('D13').Value = "'19.86"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.15%  '
$ws.Range('D14').Value = "'7.78"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.59%  '
$ws.Range('D15').Value = '3.363.41'
$ws.Range('E15').Value = '  +4.07%  '
$ws.Range('D16').Value = "'0.999"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.13%  '
$ws.Range('D17').Value = '2.903.64'
$ws.Range('E17').Value = '  +3.70%  '
$ws.Range('D18').Value = '52.281.76'
$ws.Range('E18').Value = '  +1.32%  '
$ws.Range('E19').Value = '  +1.17%  '
$ws.Range('D20').Value = "'3.31"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.54%  '
$ws.Range('D21').Value = "'14.18"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.06%  '
$ws.Range('D22').Value = '0.0₃0979'
$ws.Range('E22').Value = '  +1.15%  '
$ws.Range('D23').Value = "'70.85"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.00%  '
$ws.Range('D24').Value = "'269.93"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.08%  '
$ws.Range('E25').Value = '  +1.79%  '
$ws.Range('E26').Value = '  +8.55%  '
$ws.Range('D27').Value = "'26.81"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.03%  '
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('E29').Value = '  +17.25%  '
$ws.Range('D30').Value = "'10.63"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.01%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = "'6.59"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +8.14%  '
$ws.Range('B32').Value = 'InjectiveProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D32').Value = "'37.64"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.35%  '
$ws.Range('D33').Value = "'6.24"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +12.10%  '
$ws.Range('D34').Value = "'53.11"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.95%  '
$ws.Range('E35').Value = '  -0.52%  '
$ws.Range('D36').Value = "'1.96"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -13.08%  '
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('E38').Value = '  +6.04%  '
$ws.Range('D39').Value = "'18.95"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.07%  '
$ws.Range('E40').Value = '  +2.83%  '
$ws.Range('D41').Value = "'2.77"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +11.14%  '
$ws.Range('E42').Value = '  +1.98%  '
$ws.Range('D43').Value = "'23.03"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.69%  '
$ws.Range('D44').Value = "'2.62"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +7.34%  '
$ws.Range('D45').Value = "'119.92"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.02%  '
$ws.Range('E46').Value = '  -1.79%  '
$ws.Range('E47').Value = '  +3.94%  '
$ws.Range('D48').Value = '2.180.81'
$ws.Range('E48').Value = '  +3.59%  '
$ws.Range('D49').Value = "'0.262"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +21.41%  '
$ws.Range('D50').Value = "'0.0352"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +14.29%  '
$ws.Range('D51').Value = "'0.956"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.42%  '
